$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 496.63635
$ws.Range("I12").Value = 527.2857
$ws.Range("J12").Value = 443
$ws.Range("K12").Value = 527.2857
$ws.Range("L12").Value = 443
$ws.Range("M12").Value = -357.2857
$ws.Range("N12").Value = -783
$ws.Range("H49").Value = 2469.8333
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 2469.8333
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 7409.499899999999
$ws.Range("M49").Value = $null
$ws.Range("N49").Value = -7681.499899999999
$ws.Range("H59").Value = 973
$ws.Range("I59").Value = 1000
$ws.Range("J59").Value = 959.5
$ws.Range("K59").Value = 3000
$ws.Range("L59").Value = 2878.5
$ws.Range("M59").Value = -2443
$ws.Range("N59").Value = -3992.5
$ws.Range("H98").Value = 1313.6305
$ws.Range("I98").Value = 829.42426
$ws.Range("J98").Value = 2542.7693
$ws.Range("K98").Value = 829.42426
$ws.Range("L98").Value = 2542.7693
$ws.Range("M98").Value = 668.57574
$ws.Range("N98").Value = -5538.7693
$ws.Range("H116").Value = 3931.45
$ws.Range("I116").Value = 4254.25
$ws.Range("J116").Value = 3447.25
$ws.Range("K116").Value = 4254.25
$ws.Range("L116").Value = 3447.25
$ws.Range("M116").Value = -812.25
$ws.Range("N116").Value = -10331.25
$ws.Range("H122").Value = 1313.6305
$ws.Range("I122").Value = 829.42426
$ws.Range("J122").Value = 2542.7693
$ws.Range("K122").Value = 2488.27278
$ws.Range("L122").Value = 7628.3079
$ws.Range("M122").Value = -38.27278000000024
$ws.Range("N122").Value = -12528.3079
$ws.Range("H132").Value = 2467.319
$ws.Range("I132").Value = 1536.6
$ws.Range("K132").Value = 4609.799999999999
$ws.Range("M132").Value = -2079.799999999999
$ws.Range("H141").Value = 1634.56
$ws.Range("I141").Value = 815.6087
$ws.Range("J141").Value = 11052.5
$ws.Range("K141").Value = 2446.8261
$ws.Range("L141").Value = 33157.5
$ws.Range("M141").Value = 2733.1739
$ws.Range("N141").Value = -43517.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1021329.1
$ws.Range("I32").Value = 1228705.2
$ws.Range("K32").Value = 1228705.2
$ws.Range("M32").Value = -1228418.2
$ws.Range("H74").Value = 221659.02
$ws.Range("I74").Value = 271232.16
$ws.Range("J74").Value = 80566.234
$ws.Range("K74").Value = 271232.16
$ws.Range("L74").Value = 80566.234
$ws.Range("M74").Value = -270358.16
$ws.Range("N74").Value = -82314.234
$ws.Range("H77").Value = 221659.02
$ws.Range("I77").Value = 271232.16
$ws.Range("J77").Value = 80566.234
$ws.Range("K77").Value = 1356160.8
$ws.Range("L77").Value = 402831.17
$ws.Range("M77").Value = -1351792.8
$ws.Range("N77").Value = -411567.17
$ws.Range("H130").Value = 56582.25
$ws.Range("J130").Value = 56582.25
$ws.Range("L130").Value = 56582.25
$ws.Range("N130").Value = -66622.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 5870
$ws.Range("I38").Value = 1990
$ws.Range("J38").Value = 9750
$ws.Range("K38").Value = 1990
$ws.Range("L38").Value = 9750
$ws.Range("M38").Value = -1613
$ws.Range("N38").Value = -10504
$ws.Range("H46").Value = 5870
$ws.Range("I46").Value = 1990
$ws.Range("J46").Value = 9750
$ws.Range("K46").Value = 1990
$ws.Range("L46").Value = 9750
$ws.Range("M46").Value = -1779
$ws.Range("N46").Value = -10172
$ws.Range("H122").Value = 1725.7858
$ws.Range("I122").Value = 906.8889
$ws.Range("J122").Value = 3199.8
$ws.Range("K122").Value = 2720.6667
$ws.Range("L122").Value = 9599.400000000001
$ws.Range("M122").Value = -270.6667000000002
$ws.Range("N122").Value = -14499.4
$ws.Range("H132").Value = 1731.9811
$ws.Range("I132").Value = 909.0513
$ws.Range("J132").Value = 4024.4285
$ws.Range("K132").Value = 2727.1539
$ws.Range("L132").Value = 12073.2855
$ws.Range("M132").Value = -197.1538999999998
$ws.Range("N132").Value = -17133.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 4415.75
$ws.Range("I56").Value = 4415.75
$ws.Range("K56").Value = 4415.75
$ws.Range("M56").Value = -3885.75
$ws.Range("H121").Value = 1432.4445
$ws.Range("J121").Value = 1950.75
$ws.Range("L121").Value = 5852.25
$ws.Range("N121").Value = -8472.25
$ws.Range("H131").Value = 1468.1063
$ws.Range("J131").Value = 1611.2424
$ws.Range("L131").Value = 4833.7272
$ws.Range("N131").Value = -14913.7272

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H54").Value = 10500
$ws.Range("I54").Value = 10500
$ws.Range("J54").Value = 10500
$ws.Range("K54").Value = 10500
$ws.Range("L54").Value = 10500
$ws.Range("M54").Value = -10110
$ws.Range("N54").Value = -11280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 1038.1428
$ws.Range("I35").Value = 894.5
$ws.Range("J35").Value = 1900
$ws.Range("K35").Value = 894.5
$ws.Range("L35").Value = 1900
$ws.Range("M35").Value = -558.5
$ws.Range("N35").Value = -2572
$ws.Range("H38").Value = 8676.666999999999
$ws.Range("J38").Value = 8676.666999999999
$ws.Range("L38").Value = 8676.666999999999
$ws.Range("N38").Value = -9496.666999999999
$ws.Range("H42").Value = 8729.666999999999
$ws.Range("I42").Value = 4025
$ws.Range("J42").Value = 11082
$ws.Range("K42").Value = 4025
$ws.Range("L42").Value = 11082
$ws.Range("M42").Value = -3462
$ws.Range("N42").Value = -12208
$ws.Range("H45").Value = 4929.125
$ws.Range("I45").Value = 2847
$ws.Range("J45").Value = 6178.4
$ws.Range("K45").Value = 2847
$ws.Range("L45").Value = 6178.4
$ws.Range("M45").Value = -2440
$ws.Range("N45").Value = -6992.4
$ws.Range("H49").Value = 8729.666999999999
$ws.Range("I49").Value = 4025
$ws.Range("J49").Value = 11082
$ws.Range("K49").Value = 4025
$ws.Range("L49").Value = 11082
$ws.Range("M49").Value = -3878
$ws.Range("N49").Value = -11376
$ws.Range("H58").Value = 2000
$ws.Range("I58").Value = 2000
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 2000
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -1740
$ws.Range("N58").Value = $null
$ws.Range("H100").Value = 62505004
$ws.Range("I100").Value = 8500
$ws.Range("J100").Value = 111113390
$ws.Range("K100").Value = 8500
$ws.Range("L100").Value = 111113390
$ws.Range("M100").Value = -7959
$ws.Range("N100").Value = -111114472

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 2500
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 2500
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 2500
$ws.Range("M24").Value = $null
$ws.Range("N24").Value = -2960
$ws.Range("H49").Value = 1256926.5
$ws.Range("J49").Value = 9450
$ws.Range("L49").Value = 9450
$ws.Range("N49").Value = -9910
$ws.Range("H51").Value = 19900
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 19900
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 19900
$ws.Range("M51").Value = $null
$ws.Range("N51").Value = -20920
$ws.Range("H126").Value = 1078.8
$ws.Range("I126").Value = 740
$ws.Range("J126").Value = 2151.6667
$ws.Range("K126").Value = 2220
$ws.Range("L126").Value = 6455.000100000001
$ws.Range("M126").Value = 250
$ws.Range("N126").Value = -11395.0001
